# partial progress on #56 NOT TESTED
#
# The "designs" sheet had an extra "Lifetime" row for the Utilities process
# area duplicated in with the per-scenario blocks; row 66 (the first stray
# row, scenario block starting at row 58, Lifetime/Preprocessing) is removed
# and everything below it shifts up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("designs")

# Select the row first (so the post-delete selection lands on the row that
# slid up into the deleted row's place, matching what Excel leaves selected
# after a Delete Entire Row), then delete it, shifting cells up.
$ws.Rows.Item(66).Select() | Out-Null
$ws.Rows.Item(66).Delete() | Out-Null

# Leave the whole row selected at its new position, same as Excel does after
# a row delete.
$ws.Rows.Item(66).Select() | Out-Null
